$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$d = $ws.Range("D2")
$d.NumberFormat = "@"
$d.Value = '64.095.15'
$d.Style = "Normal"
$ws.Range("E2").Value = '  +1.39%  '

$d = $ws.Range("D3")
$d.NumberFormat = "@"
$d.Value = '3.148.15'
$d.Style = "Normal"
$ws.Range("E3").Value = '  +2.06%  '

$ws.Range("E4").Value = '  +0.10%  '

$d = $ws.Range("D5")
$d.NumberFormat = "@"
$d.Value = '590.75'
$d.Style = "Normal"
$ws.Range("E5").Value = '  +1.70%  '

$d = $ws.Range("D6")
$d.NumberFormat = "@"
$d.Value = '146.37'
$d.Style = "Normal"
$ws.Range("E6").Value = '  +1.54%  '

$d = $ws.Range("D8")
$d.NumberFormat = "@"
$d.Value = '3.137.25'
$d.Style = "Normal"
$ws.Range("E8").Value = '  +1.97%  '

$d = $ws.Range("D9")
$d.NumberFormat = "@"
$d.Value = '0.531'
$d.Style = "Normal"
$ws.Range("E9").Value = '  +0.71%  '

$ws.Range("E10").Value = '  +3.23%  '

$d = $ws.Range("D11")
$d.NumberFormat = "@"
$d.Value = '5.90'
$d.Style = "Normal"
$ws.Range("E11").Value = '  +5.33%  '

$ws.Range("E12").Value = '  +0.50%  '

$ws.Range("E13").Value = '  +1.48%  '

$d = $ws.Range("D14")
$d.NumberFormat = "@"
$d.Value = '37.17'
$d.Style = "Normal"
$ws.Range("E14").Value = '  -1.12%  '

$d = $ws.Range("D15")
$d.NumberFormat = "@"
$d.Value = '3.665.58'
$d.Style = "Normal"
$ws.Range("E15").Value = '  +1.94%  '

$ws.Range("E16").Value = '  -0.17%  '

$d = $ws.Range("D17")
$d.NumberFormat = "@"
$d.Value = '7.25'
$d.Style = "Normal"
$ws.Range("E17").Value = '  +2.39%  '

$d = $ws.Range("D18")
$d.NumberFormat = "@"
$d.Value = '63.874.54'
$d.Style = "Normal"
$ws.Range("E18").Value = '  +1.21%  '

$d = $ws.Range("D19")
$d.NumberFormat = "@"
$d.Value = '3.139.25'
$d.Style = "Normal"
$ws.Range("E19").Value = '  +1.72%  '

$d = $ws.Range("D20")
$d.NumberFormat = "@"
$d.Value = '465.99'
$d.Style = "Normal"
$ws.Range("E20").Value = '  +1.27%  '

$d = $ws.Range("D21")
$d.NumberFormat = "@"
$d.Value = '14.40'
$d.Style = "Normal"
$ws.Range("E21").Value = '  +1.42%  '

$d = $ws.Range("D22")
$d.NumberFormat = "@"
$d.Value = '0.732'
$d.Style = "Normal"
$ws.Range("E22").Value = '  +1.26%  '

$d = $ws.Range("D23")
$d.NumberFormat = "@"
$d.Value = '7.59'
$d.Style = "Normal"
$ws.Range("E23").Value = '  +2.22%  '

$d = $ws.Range("D24")
$d.NumberFormat = "@"
$d.Value = '2.39'
$d.Style = "Normal"
$ws.Range("E24").Value = '  +13.05%  '

$d = $ws.Range("D25")
$d.NumberFormat = "@"
$d.Value = '13.17'
$d.Style = "Normal"
$ws.Range("E25").Value = '  +1.89%  '

$d = $ws.Range("D26")
$d.NumberFormat = "@"
$d.Value = '80.90'
$d.Style = "Normal"
$ws.Range("E26").Value = '  -0.15%  '

$ws.Range("E27").Value = '  +0.08%  '

$d = $ws.Range("D28")
$d.NumberFormat = "@"
$d.Value = '9.85'
$d.Style = "Normal"
$ws.Range("E28").Value = '  +10.55%  '

$ws.Range("E29").Value = '  +2.08%  '

$ws.Range("E30").Value = '  +7.55%  '

$ws.Range("E31").Value = '  +0.14%  '

$ws.Range("E32").Value = '  +0.49%  '

$ws.Range("E33").Value = '  +4.73%  '

$d = $ws.Range("D34")
$d.NumberFormat = "@"
$d.Value = '27.65'
$d.Style = "Normal"
$ws.Range("E34").Value = '  +4.11%  '

$d = $ws.Range("D35")
$d.NumberFormat = "@"
$d.Value = '0.0₃0859'
$d.Style = "Normal"
$ws.Range("E35").Value = '  +2.37%  '

$ws.Range("E36").Value = '  +3.37%  '

$d = $ws.Range("D37")
$d.NumberFormat = "@"
$d.Value = '6.16'
$d.Style = "Normal"
$ws.Range("E37").Value = '  +2.95%  '

$d = $ws.Range("D38")
$d.NumberFormat = "@"
$d.Value = '2.28'
$d.Style = "Normal"
$ws.Range("E38").Value = '  -0.61%  '

$d = $ws.Range("D39")
$d.NumberFormat = "@"
$d.Value = '3.24'
$d.Style = "Normal"
$ws.Range("E39").Value = '  -2.15%  '

$d = $ws.Range("D40")
$d.NumberFormat = "@"
$d.Value = '463.31'
$d.Style = "Normal"
$ws.Range("E40").Value = '  +7.38%  '

$d = $ws.Range("D41")
$d.NumberFormat = "@"
$d.Value = '51.36'
$d.Style = "Normal"
$ws.Range("E41").Value = '  +2.36%  '

$d = $ws.Range("D42")
$d.NumberFormat = "@"
$d.Value = '9.34'
$d.Style = "Normal"
$ws.Range("E42").Value = '  +7.01%  '

$ws.Range("E43").Value = '  +8.90%  '

$ws.Range("E44").Value = '  +1.45%  '

$d = $ws.Range("D45")
$d.NumberFormat = "@"
$d.Value = '2.892.91'
$d.Style = "Normal"
$ws.Range("E45").Value = '  +1.20%  '

$d = $ws.Range("D46")
$d.NumberFormat = "@"
$d.Value = '39.97'
$d.Style = "Normal"
$ws.Range("E46").Value = '  +10.49%  '

$ws.Range("E47").Value = '  +0.14%  '

$d = $ws.Range("D48")
$d.NumberFormat = "@"
$d.Value = '132.92'
$d.Style = "Normal"
$ws.Range("E48").Value = '  +7.33%  '

$ws.Range("E50").Value = '  +0.84%  '

$ws.Range("E51").Value = '  +4.27%  '

